$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update recalculated low-flow-index statistics (correcao fim de evento)
# Columns: A = event index, D = deficit volume, F = Pexpon(<=Dobs), G = Pgama(<=Dobs)

$ws.Range("F2").Value = 0.4535277971113451
$ws.Range("G2").Value = 0.5734791886872768
$ws.Range("F3").Value = 0.1179129418069409
$ws.Range("G3").Value = 0.2853087006858848
$ws.Range("F4").Value = 0.04892478570585143
$ws.Range("G4").Value = 0.1838577298970954
$ws.Range("F5").Value = 0.03962418925596779
$ws.Range("G5").Value = 0.1655820688594732
$ws.Range("F6").Value = 0.07783742552582337
$ws.Range("G6").Value = 0.2317004224295548
$ws.Range("F7").Value = 0.451789584164324
$ws.Range("G7").Value = 0.5722734111526541
$ws.Range("F8").Value = 0.001007778286129896
$ws.Range("G8").Value = 0.02701151003266683
$ws.Range("F9").Value = 0.530018824508887
$ws.Range("G9").Value = 0.6254003625127171
$ws.Range("F10").Value = 0.378426221083617
$ws.Range("G10").Value = 0.5199529679004951
$ws.Range("F11").Value = 0.1077093146820041
$ws.Range("G11").Value = 0.2726205171527307
$ws.Range("F12").Value = 0.05419241365276442
$ws.Range("G12").Value = 0.1934466027215093
$ws.Range("F13").Value = 0.03002725819148154
$ws.Range("G13").Value = 0.1443114674944783
$ws.Range("F14").Value = 0.009236681532334645
$ws.Range("G14").Value = 0.08056955282316679
$ws.Range("F15").Value = 0.2013634510299247
$ws.Range("G15").Value = 0.3742106071737282
$ws.Range("F16").Value = 0.9474181079555035
$ws.Range("G16").Value = 0.9188942448901316
$ws.Range("D17").Value = 723.2480999999997
$ws.Range("F17").Value = 0.81647365729994
$ws.Range("G17").Value = 0.8151226585960806
$ws.Range("F18").Value = 0.8604897950091843
$ws.Range("G18").Value = 0.8467974070450799
$ws.Range("F19").Value = 0.3609808876474872
$ws.Range("G19").Value = 0.5070214096419517
$ws.Range("F20").Value = 0.8489879025846336
$ws.Range("G20").Value = 0.8383361418981622
$ws.Range("F21").Value = 0.5715572111265292
$ws.Range("G21").Value = 0.6528939738185671
$ws.Range("D22").Value = 1360.6069
$ws.Range("F22").Value = 0.9604086507376332
$ws.Range("G22").Value = 0.9321337477489449
$ws.Range("F23").Value = 0.1932713386849009
$ws.Range("G23").Value = 0.3664471160611598
$ws.Range("F24").Value = 0.859695142912881
$ws.Range("G24").Value = 0.8462078521496383
$ws.Range("F25").Value = 0.07530439785455131
$ws.Range("G25").Value = 0.2279042478599306
$ws.Range("F26").Value = 0.4922658701248102
$ws.Range("G26").Value = 0.6000282800773625
$ws.Range("F27").Value = 0.8907204146539784
$ws.Range("G27").Value = 0.8698872164845448
$ws.Range("F28").Value = 0.1533161423450484
$ws.Range("G28").Value = 0.3257362617517343
$ws.Range("F29").Value = 0.2227404065507414
$ws.Range("G29").Value = 0.3940791602247605
$ws.Range("F30").Value = 0.9828044906221561
$ws.Range("G30").Value = 0.9593250466399729
$ws.Range("F31").Value = 0.3964579111565711
$ws.Range("G31").Value = 0.5331008201872918
$ws.Range("F32").Value = 0.6452522562276861
$ws.Range("G32").Value = 0.7010601039367004
$ws.Range("F33").Value = 0.03238957470902752
$ws.Range("G33").Value = 0.1498294540412411
$ws.Range("F34").Value = 0.1306821138771109
$ws.Range("G34").Value = 0.3004740081692701
$ws.Range("D35").Value = 896.9914999999997
$ws.Range("F35").Value = 0.879185043365764
$ws.Range("G35").Value = 0.8609106099723336
$ws.Range("F36").Value = 0.4679093310140311
$ws.Range("G36").Value = 0.5834052931635964
$ws.Range("F37").Value = 0.3634216869426148
$ws.Range("G37").Value = 0.5088439645220323
$ws.Range("D38").Value = 205.312
$ws.Range("F38").Value = 0.3617528293075056
$ws.Range("G38").Value = 0.5075983035407062
$ws.Range("F39").Value = 0.3642452994200327
$ws.Range("G39").Value = 0.5094579583081865
$ws.Range("F40").Value = 0.6165322311485333
$ws.Range("G40").Value = 0.6823407821152064
$ws.Range("F41").Value = 0.006947369547073484
$ws.Range("G41").Value = 0.07000207617059133
$ws.Range("F42").Value = 0.9787383224033448
$ws.Range("G42").Value = 0.9537314912445886
$ws.Range("F43").Value = 0.2086571795210225
$ws.Range("G43").Value = 0.3810898700999379
$ws.Range("F44").Value = 0.1202672255916622
$ws.Range("G44").Value = 0.2881611314281282
$ws.Range("F45").Value = 0.042423416922967
$ws.Range("G45").Value = 0.1712876586252341
$ws.Range("F46").Value = 0.2682829346714851
$ws.Range("G46").Value = 0.4338290679515271
$ws.Range("D47").Value = 3084.823299999999
$ws.Range("F47").Value = 0.9993753498865218
$ws.Range("G47").Value = 0.9941207240915463
$ws.Range("F48").Value = 0.08188437329935869
$ws.Range("G48").Value = 0.2376428082713817
$ws.Range("F49").Value = 0.02729713181482321
$ws.Range("G49").Value = 0.1376565718724659
$ws.Range("F50").Value = 0.2698264173423163
$ws.Range("G50").Value = 0.4351239267940162
$ws.Range("F51").Value = 0.1162101971520798
$ws.Range("G51").Value = 0.2832287843870434
$ws.Range("F52").Value = 0.07706130081100657
$ws.Range("G52").Value = 0.2305437474700897
$ws.Range("F53").Value = 0.5963113107494582
$ws.Range("G53").Value = 0.6691306268091346
$ws.Range("F54").Value = 0.4339500842377791
$ws.Range("G54").Value = 0.5598173803898702
$ws.Range("F55").Value = 0.02798973924320487
$ws.Range("G55").Value = 0.1393752351845627
$ws.Range("D56").Value = 148.666
$ws.Range("F56").Value = 0.2685424323926916
$ws.Range("G56").Value = 0.434046981907281
$ws.Range("G57").Value = 0.000000003120731044221785
$ws.Range("F58").Value = 0.3292258722167992
$ws.Range("G58").Value = 0.4828777012573838
$ws.Range("F59").Value = 0.5724257799799968
$ws.Range("G59").Value = 0.6534652306011058
$ws.Range("F60").Value = 0.9857520546122556
$ws.Range("G60").Value = 0.9636870431529652
$ws.Range("A61").Value = 155
$ws.Range("F61").Value = 0.6682154361592868
$ws.Range("G61").Value = 0.7160242700099764
$ws.Range("A62").Value = 157
$ws.Range("F62").Value = 0.09117270117404355
$ws.Range("G62").Value = 0.2507645106598095
$ws.Range("A63").Value = 162
$ws.Range("F63").Value = 0.9997524188875809
$ws.Range("G63").Value = 0.9965156299081206
